$d = $word.ActiveDocument

# --- 1) "median RT TRN-RND: TRN4 - RND5 vagy residual change score" paragraph:
#     merge the "RT TRN-RND: " / "TRN4 - RND5" / " vagy " runs into one run.
$p = $d.Paragraphs(15).Range
[void]$p.Find.Execute("TRN4 – RND5", $true, $false, $false, $false, $false, $true, 1, $false, "TRN4 – RND5", 2)

# --- 2) "median RT RND-REC: REC6 - RND5 vagy residual change score" paragraph:
#     merge the "RT RND-REC: " / "REC6 - RND5" / " vagy " runs into one run.
$p = $d.Paragraphs(16).Range
[void]$p.Find.Execute("REC6 – RND5", $true, $false, $false, $false, $false, $true, 1, $false, "REC6 – RND5", 2)

# --- 3) "ACC training: TRN4 - TRN1 vagy residual change score" paragraph:
#     merge the ":" / " TRN4" / " - TRN1 vagy " runs into one run.
$p = $d.Paragraphs(17).Range
[void]$p.Find.Execute("TRN4 – TRN1", $true, $false, $false, $false, $false, $true, 1, $false, "TRN4 – TRN1", 2)

# --- 4) "ACC TRN-RND: RND5 - TRN4 vagy residual change score" paragraph:
#     merge the "ACC TRN-RND: " / "RND5 - TRN4" / " vagy " runs into one run.
$p = $d.Paragraphs(18).Range
[void]$p.Find.Execute("RND5 – TRN4", $true, $false, $false, $false, $false, $true, 1, $false, "RND5 – TRN4", 2)

# --- 5) "ACC RND-REC: RND5 - REC6 vagy residual change score" paragraph:
#     merge the "ACC RND-REC: " / "RND5 - REC6" / " vagy " runs into one run.
$p = $d.Paragraphs(19).Range
[void]$p.Find.Execute("RND5 – REC6", $true, $false, $false, $false, $false, $true, 1, $false, "RND5 – REC6", 2)

# --- 6) Collapse the whole "1. Előszűrés ... / 2. Online indexeknél ... / 3. Offline
#     feladaton belül: ... / 4. ?" block (paragraphs 25-31) down to a single
#     paragraph reading "Lásd AGL" (keeping the original last paragraph's mark /
#     pPr, which carries the bookmark "_GoBack").
$startPara = $d.Paragraphs(25)
$endPara = $d.Paragraphs(31)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
$delRange.Delete()

$last = $d.Paragraphs(25).Range
[void]$last.MoveEnd(1, -1)
$last.Text = "Lásd AGL"
